$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41 (pushes existing rows 41+ down to 42+,
# auto-adjusting formulas / merged cells / row count).
$ws.Rows("41:41").Insert()

# The new row 41 should look like the "Nama Pertandingan" detail row that is
# now at row 40 (same formatting). Copy formats only from row 40 into row 41.
$ws.Range("A40:G40").Copy()
$ws.Range("A41:G41").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the values for the newly inserted detail row.
$ws.Range("C41").Value = "Nama Pertandingan"
$ws.Range("D41").Value = 0.0
$ws.Range("E41").Value = 0.0

# Match row 40's row height exactly.
$ws.Rows("41:41").RowHeight = 12.75

# D40 picks up the same fill/border/alignment style as C40 (s=38).
$ws.Range("C40").Copy()
$ws.Range("D40").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D40").Value = 0

# The subtotal formula (now on row 42) should include the new row, inserted
# just before the D40-E40 term.
$ws.Range("F42").Formula = "=D34-E34+D35-E35+D36-E36+D37-E37+D38-E38+D39-E39+D41-E41+D40-E40"
